# Select the "Coupling Parameters" sheet (already active) and add the new
# "realistic_candidate_capacities_for_future" parameter row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Write the description (column C) first, then the key name (column A), then
# the boolean value (column B) -- this mirrors the order the new shared
# strings were originally authored in.
$ws.Range("C17").Value = "If this is true, the real capacity of the candidate power plants is considered for the future investments"
$ws.Range("A17").Value = "realistic_candidate_capacities_for_future"
$ws.Range("B17").Value = $true

# Move the selection, matching the cursor position left behind by the editor.
$ws.Range("C37").Select() | Out-Null
